# Update cryptos list - GitHub Actions data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Bitcoin ---
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.333.44"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +2.50%  "

# --- Row 3: Ethereum ---
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.343.43"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +3.48%  "

# --- Row 5: Solana ---
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "193.25"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +5.55%  "

# --- Row 6: BNB ---
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "591.18"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +2.18%  "

# --- Row 7: USDC ---
$ws.Range("E7").Value = "  -0.02%  "

# --- Row 8: XRP ---
$ws.Range("E8").Value = "  +0.82%  "

# --- Row 9: Dogecoin ---
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.134"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +3.42%  "

# --- Row 10: Toncoin ---
$ws.Range("E10").Value = "  +2.24%  "

# --- Row 11: Cardano ---
$ws.Range("E11").Value = "  +2.06%  "

# --- Row 12: Wrapped liquid staked Ether 2.0 ---
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "3.924.60"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +3.59%  "

# --- Row 13: TRON ---
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.139"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.97%  "

# --- Row 14: Avalanche ---
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.19"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +2.66%  "

# --- Row 15: Wrapped BTC ---
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "69.317.01"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +2.41%  "

# --- Row 16: Shiba Inu ---
$ws.Range("E16").Value = "  +1.74%  "

# --- Row 17: Wrapped Ether ---
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.342.78"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +2.83%  "

# --- Row 18: Polkadot ---
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.82"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +1.42%  "

# --- Row 19: Chainlink ---
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.75"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +2.53%  "

# --- Row 20: Bitcoin Cash ---
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "432.43"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +8.72%  "

# --- Row 21: Uniswap ---
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.75"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +2.92%  "

# --- Row 22: Litecoin ---
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.12"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +2.98%  "

# --- Row 23: Dai ---
$ws.Range("E23").Value = "  +0.05%  "

# --- Row 24: Wrapped eETH ---
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.494.89"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +3.44%  "

# --- Row 25: Polygon ---
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.517"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.97%  "

# --- Row 26: PEPE ---
$ws.Range("E26").Value = "  +3.55%  "

# --- Row 27: Kaspa ---
$ws.Range("E27").Value = "  +3.65%  "

# --- Row 28: Internet Computer (DFINITY) ---
$ws.Range("E28").Value = "  +0.75%  "

# --- Row 29: Binance-Peg BSC-USD ---
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.03"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +2.81%  "

# --- Row 30: PancakeSwap ---
$ws.Range("E30").Value = "  +2.39%  "

# --- Row 31: Ethereum Classic ---
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "23.03"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +1.89%  "

# --- Row 32: NEAR Protocol ---
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.59"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.82%  "

# --- Row 33: Fetch.AI ---
$ws.Range("E33").Value = "  +2.22%  "

# --- Row 34: Aptos ---
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.99"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +0.34%  "

# --- Row 35: USDe ---
$ws.Range("E35").Value = "  +0.03%  "

# --- Rows 36-37: Monero and ImmutableX swapped positions (ranking changed) ---
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.52"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +3.36%  "

$ws.Range("B37").Value = "Monero"
$ws.Range("C37").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "164.68"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +2.12%  "

# --- Row 38: Stacks ---
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.92"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +2.39%  "

# --- Row 39: EnergySwap ---
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "27.07"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +2.90%  "

# --- Row 40: Mantle ---
$ws.Range("E40").Value = "  +0.67%  "

# --- Row 41: Filecoin ---
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.57"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.54%  "

# --- Rows 42-43: RenderToken and Maker swapped positions (ranking changed) ---
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.754.60"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +6.04%  "

$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.49"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +0.79%  "

# --- Row 44: dogwifhat ---
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.52"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +3.14%  "

# --- Row 45: OKB ---
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "41.20"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.26%  "

# --- Row 46: Hedera ---
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0687"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.88%  "

# --- Row 47: Bittensor ---
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "344.19"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +3.04%  "

# --- Row 48: Injective Protocol ---
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "25.26"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +2.10%  "

# --- Row 49: VeChain ---
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0283"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +2.48%  "

# --- Row 50: Arweave ---
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "32.59"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +5.89%  "

# --- Row 51: ONDO ---
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.01"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +3.86%  "

